# Update the "Förändrad" (Changed) date column for rows 2-5
# from 2023-09-15 (serial 45184) to 2023-09-16 (serial 45185).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2").Value = 45185
$ws.Range("C3").Value = 45185
$ws.Range("C4").Value = 45185
$ws.Range("C5").Value = 45185
